$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "./model_output/2025-08-18-11-18-17-None"
$ws.Range("B7").Value = 0.6428571428571428

$ws.Range("A8").Value = "./model_output/2025-08-18-13-54-24-None"
$ws.Range("B8").Value = 0.6785714285714286

$ws.Range("A9").Value = "./model_output/2025-08-18-14-14-26-None"
$ws.Range("B9").Value = 0.5785714285714286
